$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '56.967.77'
$ws.Range("E2").Value = '  +0.08%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.412.29'
$ws.Range("E3").Value = '  -3.60%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '485.90'
$ws.Range("E5").Value = '  -2.05%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '153.18'
$ws.Range("E6").Value = '  -0.21%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.997'
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.608'
$ws.Range("E8").Value = '  +17.87%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.432.15'
$ws.Range("E9").Value = '  -3.39%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0995'
$ws.Range("E10").Value = '  +0.25%  '
$ws.Range("E11").Value = '  -2.76%  '
$ws.Range("E12").Value = '  -0.55%  '
$ws.Range("E13").Value = '  +1.25%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.836.87'
$ws.Range("E14").Value = '  -3.42%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '57.064.74'
$ws.Range("E15").Value = '  +0.16%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.74'
$ws.Range("E16").Value = '  -3.14%  '
$ws.Range("E17").Value = '  -3.20%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.428.69'
$ws.Range("E18").Value = '  -3.34%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.72'
$ws.Range("E19").Value = '  +3.75%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '324.19'
$ws.Range("E20").Value = '  -0.09%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.95'
$ws.Range("E21").Value = '  -3.88%  '
$ws.Range("E22").Value = '  -0.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.88'
$ws.Range("E23").Value = '  -1.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '58.16'
$ws.Range("E24").Value = '  -1.60%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.408'
$ws.Range("E25").Value = '  -1.14%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.41%  '
$ws.Range("E27").Value = '  -3.81%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.525.51'
$ws.Range("E28").Value = '  -3.26%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.23'
$ws.Range("E29").Value = '  -6.27%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0781'
$ws.Range("E30").Value = '  -4.74%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.999'
$ws.Range("E31").Value = '  +0.08%  '
$ws.Range("E32").Value = '  +0.81%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '148.48'
$ws.Range("E33").Value = '  -1.99%  '
$ws.Range("E34").Value = '  -1.05%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.32'
$ws.Range("E35").Value = '  +1.00%  '
$ws.Range("E36").Value = '  -2.47%  '
$ws.Range("E37").Value = '  -2.76%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.845'
$ws.Range("E38").Value = '  -4.02%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '34.10'
$ws.Range("E39").Value = '  -0.72%  '
$ws.Range("E40").Value = '  +8.99%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.51'
$ws.Range("E41").Value = '  -0.41%  '
$ws.Range("E42").Value = '  -1.86%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.997'
$ws.Range("E43").Value = '  +0.22%  '
$ws.Range("B45").Value = 'Hedera'
$ws.Range("C45").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0531'
$ws.Range("E45").Value = '  -6.49%  '
$ws.Range("B46").Value = 'Bittensor'
$ws.Range("C46").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '264.23'
$ws.Range("E46").Value = '  -1.85%  '
$ws.Range("B47").Value = 'WhiteBITCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.22'
$ws.Range("E47").Value = '  -0.07%  '
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.68'
$ws.Range("E48").Value = '  -5.26%  '
$ws.Range("E49").Value = '  -1.34%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.41'
$ws.Range("E50").Value = '  -3.10%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.862.37'
$ws.Range("E51").Value = '  -2.46%  '
